# Training / Exam dashboard refresh:
#  - header + title rows get a bold white font (drop the old size-14 title font)
#  - "PERIOD TO EXPIRE" / "LAST UPDATE" columns roll forward to the 16-Sep-2025 run
#  - Exam Dashboard comments column widened and re-worded
#  - Exam Dashboard "OK" comments reworded to "date is valid"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# --- styles.xml: header (s=2) + title (s=1) rows both become bold + white -------
$ws1.Range("A2:K2").Font.Bold = $true
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Bold = $true
$ws2.Range("A2:G2").Font.Color = 16777215

$ws1.Range("A1").Font.Bold = $true
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215

# --- Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) ---
$periods = @{
    3  = 470
    4  = 243
    5  = 526
    6  = 363
    7  = 244
    8  = 525
    9  = 392
    10 = 503
    11 = 489
    12 = 710
    13 = 527
    14 = 255
    15 = 399
    16 = 706
    17 = 489
    18 = -23
    19 = -103
    20 = -180
    21 = -36
    22 = -36
    23 = 170
    24 = 155
    25 = 278
    26 = 323
    27 = 348
}

foreach ($row in $periods.Keys) {
    $ws1.Cells.Item($row, 8).Value = $periods[$row]
    # leading apostrophe forces literal text so "16-Sep-2025" isn't parsed into a date serial
    $ws1.Cells.Item($row, 9).Value = "'16-Sep-2025"
}

# --- Exam Dashboard: widen the COMMENTS column ----------------------------------
$ws2.Columns("E").ColumnWidth = 14.15

# --- Exam Dashboard: reword the per-row comments --------------------------------
for ($row = 3; $row -le 13; $row++) {
    $ws2.Cells.Item($row, 5).Value = "date is valid"
}
